$d = $word.ActiveDocument

# --- Update the date line ---
$d.Content.Find.Execute("2025-06-28 Saturday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2025-06-29 Sunday", 2) | Out-Null

# --- Update the multiplication table cells ---
# The table has 20 rows x 5 columns; only rows 1, 5, 10, 15, 20 (1-based)
# contain data, with blank filler rows in between.
$t = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; New="853×8="},
    @{Row=1;  Col=2; New="638×7="},
    @{Row=1;  Col=3; New="364×8="},
    @{Row=1;  Col=4; New="726×3="},
    @{Row=1;  Col=5; New="418×6="},

    @{Row=5;  Col=1; New="864×4="},
    @{Row=5;  Col=2; New="508×5="},
    @{Row=5;  Col=3; New="258×5="},
    @{Row=5;  Col=4; New="975×5="},
    @{Row=5;  Col=5; New="570×5="},

    @{Row=10; Col=1; New="655×9="},
    @{Row=10; Col=2; New="615×5="},
    @{Row=10; Col=3; New="219×8="},
    @{Row=10; Col=4; New="289×6="},
    @{Row=10; Col=5; New="777×8="},

    @{Row=15; Col=1; New="684×9="},
    @{Row=15; Col=2; New="773×8="},
    @{Row=15; Col=3; New="121×9="},
    @{Row=15; Col=4; New="756×8="},
    @{Row=15; Col=5; New="435×2="},

    @{Row=20; Col=1; New="809×2="},
    @{Row=20; Col=2; New="436×4="},
    @{Row=20; Col=3; New="525×2="},
    @{Row=20; Col=4; New="977×6="},
    @{Row=20; Col=5; New="468×4="}
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.New
}
